# Add a new teammate / row (row 3) with their repo link and the six
# team-member names, matching the "add names of teammate" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Team member names (B3:G3)
$ws.Range("B3").Value = "عبدالرحمن حسين أحمد"
$ws.Range("C3").Value = "محمد الألفي محمد شريف"
$ws.Range("D3").Value = "محمد بهاء محمد عطيه"
$ws.Range("E3").Value = "أحمد غنيمي حلمي غنيمي"
$ws.Range("F3").Value = "سعيد مجدي سعيد محمدي"
$ws.Range("G3").Value = "حازم خالد منصور بيومي"

# Repo link cell (A3) - hyperlink whose display text is the URL itself,
# matching the existing A2 pattern.
$repoUrl = "https://github.com/Abdo3882/Open-source-pro"
$ws.Range("A3").Value = $repoUrl
$ws.Hyperlinks.Add($ws.Range("A3"), $repoUrl, [Type]::Missing, [Type]::Missing, $repoUrl)

# Move the active selection to A5, matching the saved view state.
$ws.Range("A5").Select()

Write-Host "done"
